$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: one before the old "Saldo Simpanan" column (C),
# and one before the old "Saldo Pinjaman" column (which becomes E after the
# first insert). This shifts the existing data right instead of overwriting
# any existing cells in place.
$ws.Columns("C:C").Insert()
$ws.Columns("E:E").Insert()

# Fix up the header style of the very first new column: Insert() copied the
# format from the column to its left (matching "Nomor Anggota"'s style),
# but it should match its right neighbour ("Saldo Simpanan"'s header style).
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "Keterangan" header cells.
$ws.Range("C1").Value2 = "Keterangan"
$ws.Range("E1").Value2 = "Keterangan"

# New "IURAN BULAN JUNI 2025" note cells.
$ws.Range("C2").Value2 = "IURAN BULAN JUNI 2025"
$ws.Range("E2").Value2 = "IURAN BULAN JUNI 2025"

# Match the column widths used in the updated template.
$ws.Columns("C:C").ColumnWidth = 21
$ws.Columns("D:D").ColumnWidth = 13
$ws.Columns("E:E").ColumnWidth = 21
$ws.Columns("F:F").ColumnWidth = 12.666666666666666

# Move the active selection like the saved workbook.
$ws.Range("E6").Select()
